$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(42).Insert()

$ws.Range("A42").Value = 8
$ws.Range("B42").Value = "Terminal La Palmera de La Serena"
$ws.Range("C42").Value = "Coquimbo"
$ws.Range("D42").Value = 44447
$ws.Range("E42").Value = 4
$ws.Range("F42").Value = 100112021
$ws.Range("G42").Value = "Ají"
$ws.Range("H42").Value = "Inferno"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 600
$ws.Range("K42").Value = 42000
$ws.Range("L42").Value = 43000
$ws.Range("M42").Value = 42500
$ws.Range("N42").Value = "$/caja 12 kilos"
$ws.Range("O42").Value = "Región de Arica y Parinacota"
$ws.Range("P42").Value = 3542
$ws.Range("Q42").Value = 12
$ws.Range("R42").Value = "Hortaliza"
